$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (post-edit) table of players, positions and teams, rows 2-19.
# Each player keeps their own Position/Team; the row order was changed
# and "Luguentz Dort" was replaced by "Robert Williams III".
$data = @(
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("Donovan Clingan", "C", "Portland Trail Blazers"),
    @("Robert Williams III", "C", "Portland Trail Blazers"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("Goga Bitadze", "C", "Orlando Magic"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers")
)

# Write column by column (all names, then all positions, then all teams)
# so that new shared strings are appended to the table in the same order
# a human re-typing the sheet column-by-column would produce.
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
